# backup before use the imp_th_process
#
# Adds two new rows ("ny_start" / "ny_end") inside the mesh parameter block
# (directly below the existing "ny" row), clears the stray AssmShellThick
# value, and refreshes the active view/selection - matching the prep work
# done before wiring up imp_th_process.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlEdgeLeft        = 7
$xlEdgeTop         = 8
$xlEdgeBottom      = 9
$xlEdgeRight       = 10
$xlContinuous      = 1
$xlThin            = 2
$xlNone            = -4142
$xlCenter          = -4108
$xlLeft            = -4131

function Set-BoxBorder($addr) {
    $ws.Range($addr).BorderAround($xlContinuous, $xlThin)
}

function Set-SideBorders($addr, [bool]$left, [bool]$top, [bool]$right, [bool]$bottom) {
    $r = $ws.Range($addr)
    if ($left) {
        $r.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
        $r.Borders.Item($xlEdgeLeft).Weight = $xlThin
    } else {
        $r.Borders.Item($xlEdgeLeft).LineStyle = $xlNone
    }
    if ($right) {
        $r.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
        $r.Borders.Item($xlEdgeRight).Weight = $xlThin
    } else {
        $r.Borders.Item($xlEdgeRight).LineStyle = $xlNone
    }
    if ($top) {
        $r.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
        $r.Borders.Item($xlEdgeTop).Weight = $xlThin
    } else {
        $r.Borders.Item($xlEdgeTop).LineStyle = $xlNone
    }
    if ($bottom) {
        $r.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
        $r.Borders.Item($xlEdgeBottom).Weight = $xlThin
    } else {
        $r.Borders.Item($xlEdgeBottom).LineStyle = $xlNone
    }
}

# --- Insert two new rows right after row 16 (the existing "ny" row), which
#     pushes the hydraulic/boundary/init/confactor blocks down by two rows.
#     Excel carries the row-16 formatting onto the freshly inserted rows,
#     and auto-shifts the three lower merged groups (G18:G21, G22:G24,
#     G25:G26) without any help - only the G13:G16 "mesh" merge needs to be
#     grown by hand afterwards.
$ws.Rows("17:18").Insert()

# Row 16 ("ny") no longer carries the unit/remark text - that now lives on
# the new "ny_end" row below it.
$ws.Range("J16").ClearContents()
$ws.Range("K16").ClearContents()

# --- New row 17: ny_start (middle piece of the G13:G18 merge band)
Set-SideBorders "G17" $true $false $true $false
$ws.Range("G17").HorizontalAlignment = $xlCenter
$ws.Range("G17").VerticalAlignment = $xlCenter

$ws.Range("H17").Value = "ny_start"
$ws.Range("I17").Value = 1
$ws.Range("J17").ClearContents()
$ws.Range("K17").ClearContents()

Set-BoxBorder "H17"
Set-BoxBorder "I17"
Set-BoxBorder "J17"
Set-BoxBorder "K17"
$ws.Range("H17:K17").VerticalAlignment = $xlCenter
$ws.Range("I17").Locked = $false
$ws.Range("I17").HorizontalAlignment = $xlLeft

# --- New row 18: ny_end (bottom piece of the G13:G18 merge band)
Set-SideBorders "G18" $true $false $true $true
$ws.Range("G18").HorizontalAlignment = $xlCenter
$ws.Range("G18").VerticalAlignment = $xlCenter

$ws.Range("H18").Value = "ny_end"
$ws.Range("I18").Value = 10
$ws.Range("J18").Value = "个"
$ws.Range("K18").Value = "轴向控制体数"

Set-SideBorders "H18" $true $false $true $false
$ws.Range("H18").VerticalAlignment = $xlCenter

Set-SideBorders "I18" $true $false $true $false
$ws.Range("I18").VerticalAlignment = $xlCenter
$ws.Range("I18").HorizontalAlignment = $xlLeft
$ws.Range("I18").Locked = $false

Set-BoxBorder "J18"
Set-BoxBorder "K18"
$ws.Range("J18:K18").VerticalAlignment = $xlCenter

# --- Grow the "mesh" group's merged label cell to cover the two new rows.
$ws.Range("G13:G18").Merge()

# --- Clear the stray AssmShellThick value (I8) - now blank in the template.
$ws.Range("I8").ClearContents()

# --- Sheet view: scrolled so row 4 is at the top, with E10 selected.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("E10").Select()
